$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Full weapon data set (weaponID, damage, reloadSpeed, magazine, useTime, bulletSpeed)
$weaponData = @(
    @("RifleS1",   12, 1, 40, 0.1, 10),
    @("RifleA1",    8, 2, 40, 0.2, 20),
    @("RifleB1",    6, 3, 40, 0.3, 30),
    @("RifleC1",    4, 4, 40, 0.4, 40),
    @("PistolS1",  30, 1, 12, 1,   10),
    @("PistolA1",  25, 2, 12, 1,   20),
    @("PistolB1",  20, 3, 12, 1,   30),
    @("PistolC1",  15, 4, 12, 1,   40),
    @("SnipeS1",   12, 1, 3,  1.5, 10),
    @("SnipeA1",    8, 2, 3,  1.5, 20),
    @("SnipeB1",    6, 3, 3,  1.5, 30),
    @("SnipeC1",    4, 4, 3,  1.5, 40),
    @("ShotgunS1", 12, 1, 4,  1,   10),
    @("ShotgunA1",  8, 2, 2,  0.1, 20),
    @("ShotgunB1",  6, 3, 4,  1,   30),
    @("ShotgunC1",  4, 4, 4,  1,   40)
)

# weaponID (column A) was originally typed in the order RifleA1, RifleS1, ... -
# enter row 3's id before row 2's so the shared-string table gets the same order.
$ws1.Cells.Item(3, 1).Value = "RifleA1"
$ws1.Cells.Item(2, 1).Value = "RifleS1"

$row = 2
foreach ($w in $weaponData) {
    $ws1.Cells.Item($row, 1).Value = $w[0]
    $ws1.Cells.Item($row, 2).Value = $w[1]
    $ws1.Cells.Item($row, 3).Value = $w[2]
    $ws1.Cells.Item($row, 4).Value = $w[3]
    $ws1.Cells.Item($row, 5).Value = $w[4]
    $ws1.Cells.Item($row, 6).Value = $w[5]
    $row++
}

# Widen the weaponID column to fit the new ids
$ws1.Columns.Item(1).ColumnWidth = 10

# Make WeaponDB the active sheet/tab with H10 selected, EnemyDB no longer selected
$ws1.Activate()
[void]$ws1.Range("H10").Select()
